# Repull data, push all data, mean calculation
# Update the "dSF" (column F) values to reflect the re-pulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new F-column (dSF) value
$updates = @{
    3  = 6
    4  = 4
    8  = 1
    11 = 3
    23 = -3
    29 = 2
    33 = 2
    36 = 2
    42 = 0
    44 = 3
    45 = 0
    46 = 0
    49 = 0
    51 = -2
    53 = 6
    56 = 1
    58 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
